# Update data: append the daily COVID-19 figures for 2020-06-08 (serial 43969)
# as a new row at the bottom of the "Tabela1" table, extending the table/
# autofilter range from A1:J68 to A1:J69.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Grow the table by one row (this also extends ref/autoFilter and dimension).
$newRow = $tbl.ListRows.Add()

# Copy the formatting of two rows up (row 67) so the new last row keeps the
# same banding/border look the sheet already uses for its final rows.
$ws.Range("A67:J67").Copy() | Out-Null
$ws.Range("A69:J69").PasteSpecial(-4122) | Out-Null

# Write the new day's values.
$r = 69
$ws.Cells.Item($r, 1).Value  = 43969   # Date
$ws.Cells.Item($r, 2).Value  = 70970   # Tested (all)
$ws.Cells.Item($r, 3).Value  = 1128    # Tested (daily)
$ws.Cells.Item($r, 4).Value  = 1467    # Positive (all)
$ws.Cells.Item($r, 5).Value  = 1       # Positive (daily)
$ws.Cells.Item($r, 6).Value  = 24      # All hospitalized on certain day
$ws.Cells.Item($r, 7).Value  = 5       # All persons in intensive care
$ws.Cells.Item($r, 8).Value  = 1       # Discharged
$ws.Cells.Item($r, 9).Value  = 104     # Deaths (all)
$ws.Cells.Item($r, 10).Value = 0       # Deaths (daily)

# Match Excel's behaviour of moving the active selection to the new last row.
$ws.Range("A69:J69").Select() | Out-Null
